$wb = $excel.ActiveWorkbook

# Sheet "fig1_case_maps": update cumulative case counts for 2020-04-06
$wsCaseMaps = $wb.Worksheets.Item("fig1_case_maps")
$wsCaseMaps.Range("D484").Value = 16524
$wsCaseMaps.Range("D506").Value = 318
$wsCaseMaps.Range("D526").Value = 2878

# Sheet "fig1_epi_timeseries": update USA cases/deaths timeseries (adding "res" to full run)
$wsEpi = $wb.Worksheets.Item("fig1_epi_timeseries")
$wsEpi.Range("B2").Value = 81
$wsEpi.Range("C2").Value = 9
$wsEpi.Range("B3").Value = 116
$wsEpi.Range("C3").Value = 13
$wsEpi.Range("B4").Value = 182
$wsEpi.Range("C4").Value = 14
$wsEpi.Range("B5").Value = 231
$wsEpi.Range("C5").Value = 19
$wsEpi.Range("B6").Value = 356
$wsEpi.Range("C6").Value = 22
$wsEpi.Range("B7").Value = 471
$wsEpi.Range("C7").Value = 25
$wsEpi.Range("B8").Value = 733
$wsEpi.Range("C8").Value = 31
$wsEpi.Range("B9").Value = 1018
$wsEpi.Range("C9").Value = 34
$wsEpi.Range("B10").Value = 1361
$wsEpi.Range("C10").Value = 40
$wsEpi.Range("B11").Value = 1778
$wsEpi.Range("C11").Value = 44
$wsEpi.Range("B12").Value = 2321
$wsEpi.Range("C12").Value = 51
$wsEpi.Range("B13").Value = 2952
$wsEpi.Range("C13").Value = 60
$wsEpi.Range("B14").Value = 3733
$wsEpi.Range("C14").Value = 71
$wsEpi.Range("B15").Value = 4710
$wsEpi.Range("C15").Value = 99
$wsEpi.Range("B16").Value = 6281
$wsEpi.Range("C16").Value = 115
$wsEpi.Range("B17").Value = 9199
$wsEpi.Range("C17").Value = 150
$wsEpi.Range("B18").Value = 14221
$wsEpi.Range("C18").Value = 209
$wsEpi.Range("B19").Value = 19132
$wsEpi.Range("C19").Value = 251
$wsEpi.Range("B20").Value = 26459
$wsEpi.Range("C20").Value = 328
$wsEpi.Range("B21").Value = 33747
$wsEpi.Range("B22").Value = 43884
$wsEpi.Range("C22").Value = 549
$wsEpi.Range("B23").Value = 53877
$wsEpi.Range("C23").Value = 739
$wsEpi.Range("B24").Value = 68561
$wsEpi.Range("C24").Value = 1004
$wsEpi.Range("B25").Value = 84969
$wsEpi.Range("C25").Value = 1268
$wsEpi.Range("B26").Value = 102598
$wsEpi.Range("C26").Value = 1615
$wsEpi.Range("B27").Value = 122971
$wsEpi.Range("C27").Value = 2112
$wsEpi.Range("B28").Value = 141740
$wsEpi.Range("C28").Value = 2461
$wsEpi.Range("B29").Value = 163692
$wsEpi.Range("C29").Value = 3005
$wsEpi.Range("B30").Value = 186736
$wsEpi.Range("C30").Value = 3849
$wsEpi.Range("B31").Value = 213121
$wsEpi.Range("C31").Value = 4788
$wsEpi.Range("B32").Value = 241843
$wsEpi.Range("C32").Value = 5877
$wsEpi.Range("B33").Value = 277047
$wsEpi.Range("C33").Value = 7054
$wsEpi.Range("B34").Value = 310725
$wsEpi.Range("C34").Value = 8234
$wsEpi.Range("B35").Value = 335443
$wsEpi.Range("C35").Value = 9539
$wsEpi.Range("B36").Value = 365307
